# Two new customer-review rows were added to the feedback export
# ("avaliacoes_garantia"): one lands at row 10 and another at row 15 in
# the final sheet, pushing all the rows that used to be 10-16 down by
# one (and the ones from the second insertion point down by one more).
# Net effect: sheet grows from A1:D16 to A1:D18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first new row at position 10 -------------------------
# (shifts old rows 10..16 down to 11..17)
$ws.Rows("10:10").Insert()

$ws.Cells.Item(10, 1).Value = 5
$ws.Cells.Item(10, 2).Value = ""
$ws.Cells.Item(10, 3).Value = 45919.59230513889
$ws.Cells.Item(10, 4).Value = "MDZiY2UyZTQtNDM3Yy00MmY5LTk2NDktMzcwNzU4YWZjZWZlOjU3MDE2"

# --- Insert the second new row at position 15 -------------------------
# (after the first insert, the old rows 10..16 now sit at 11..17; this
# shifts the old-13..16 block, now at rows 14..17, down to 15..18)
$ws.Rows("15:15").Insert()

$ws.Cells.Item(15, 1).Value = 5
$ws.Cells.Item(15, 2).Value = "Parece que vocês estão dando um rumo melhor no atendimento "
$ws.Cells.Item(15, 3).Value = 45919.58513085648
$ws.Cells.Item(15, 4).Value = "NzA2MGY5NmYtZmUwZi00NTc2LTlmMDAtZjMzM2NjYzJmYWI1OjU3MDE2"
